# doc new OutputValues option
# Add a new row to the "Description" sheet documenting the OutputValues
# trait-file option (Name / Description / Value columns), matching the
# existing table's layout and formatting, then refresh the row heights
# that Excel recomputed for the (now slightly re-wrapped) existing rows,
# and leave the selection on the new last cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

# --- New row 16: copy the formatting of the last existing row (15) ---
$ws.Range("A15:C15").Copy()
$null = $ws.Range("A16:C16").PasteSpecial(-4122)

# Fill in cell values in this order so new shared-string entries land in
# the same order Excel produced them: Name, Value, then Description.
$ws.Range("A16").Value = "OutputValues"
$ws.Range("C16").Value = "TRUE/FALSE"
$ws.Range("B16").Value = "Should allele values for this gene be written to output? Ignored if OutputGeneValues in GeneticsFile is set to FALSE. "

$ws.Rows.Item(16).RowHeight = 45

# --- Row height refresh on existing rows (re-wrap after edit) ---
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 90
$ws.Rows.Item(5).RowHeight = 90
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 105
$ws.Rows.Item(9).RowHeight = 60
$ws.Rows.Item(10).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 105
$ws.Rows.Item(12).RowHeight = 90
$ws.Rows.Item(13).RowHeight = 150
$ws.Rows.Item(14).RowHeight = 105

# Column B widened (and no longer "best fit") to comfortably fit the
# new option's longer description text.
$ws.Columns.Item(2).ColumnWidth = 49.166666666666664

# Move the selection to the new last cell, as in the saved file.
$null = $ws.Range("C16").Select()
